# Apply scheduled-runner updates to Seraph_Profits sheets.
# Each sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) has a handful of
# profit-calculation rows whose cached H/I/J/K/L/M/N numbers need to be
# refreshed to match the latest recipe/market data.

$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4998.5
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H62").Value = 4091.1428
$ws.Range("I62").Value = 1919.5555
$ws.Range("K62").Value = 1919.5555
$ws.Range("M62").Value = -1295.5555
$ws.Range("H64").Value = 6500
$ws.Range("I64").Value = 6000
$ws.Range("J64").Value = 8000
$ws.Range("K64").Value = 6000
$ws.Range("L64").Value = 8000
$ws.Range("M64").Value = -5752
$ws.Range("N64").Value = -8496
$ws.Range("H65").Value = 4091.1428
$ws.Range("I65").Value = 1919.5555
$ws.Range("K65").Value = 9597.7775
$ws.Range("M65").Value = -6477.7775
$ws.Range("H67").Value = 6500
$ws.Range("I67").Value = 6000
$ws.Range("J67").Value = 8000
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 8000
$ws.Range("M67").Value = -5142
$ws.Range("N67").Value = -9716
$ws.Range("H70").Value = 97271.375
$ws.Range("I70").Value = 1166.3334
$ws.Range("K70").Value = 3499.0002
$ws.Range("M70").Value = -3229.0002
$ws.Range("H73").Value = 97271.375
$ws.Range("I73").Value = 1166.3334
$ws.Range("K73").Value = 3499.0002
$ws.Range("M73").Value = -2563.0002
$ws.Range("H137").Value = 2289
$ws.Range("I137").Value = 1116.25
$ws.Range("K137").Value = 3348.75
$ws.Range("M137").Value = -798.75

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5210.6665
$ws.Range("I74").Value = 1581
$ws.Range("K74").Value = 1581
$ws.Range("M74").Value = -707
$ws.Range("H77").Value = 5210.6665
$ws.Range("I77").Value = 1581
$ws.Range("K77").Value = 7905
$ws.Range("M77").Value = -3537
$ws.Range("H97").Value = 740.44446
$ws.Range("I97").Value = 787.8
$ws.Range("J97").Value = 681.25
$ws.Range("K97").Value = 787.8
$ws.Range("L97").Value = 681.25
$ws.Range("M97").Value = -291.8
$ws.Range("N97").Value = -1673.25
$ws.Range("H102").Value = 1667.091
$ws.Range("I102").Value = 941.2857
$ws.Range("K102").Value = 941.2857
$ws.Range("M102").Value = 680.7143
$ws.Range("H110").Value = 4834.5
$ws.Range("I110").Value = 4834.5
$ws.Range("K110").Value = 4834.5
$ws.Range("M110").Value = -2789.5
$ws.Range("H122").Value = 913089.4399999999
$ws.Range("I122").Value = 1669830.9
$ws.Range("K122").Value = 5009492.699999999
$ws.Range("M122").Value = -5007042.699999999
$ws.Range("H132").Value = 1654.5834
$ws.Range("I132").Value = 1577.9697
$ws.Range("K132").Value = 4733.909100000001
$ws.Range("M132").Value = -2203.909100000001

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1036.375
$ws.Range("I107").Value = 882.6667
$ws.Range("J107").Value = 1497.5
$ws.Range("K107").Value = 882.6667
$ws.Range("L107").Value = 1497.5
$ws.Range("M107").Value = 1037.3333
$ws.Range("N107").Value = -5337.5

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5205.1875
$ws.Range("I31").Value = 3367.375
$ws.Range("K31").Value = 3367.375
$ws.Range("M31").Value = -3072.375
$ws.Range("H34").Value = 5205.1875
$ws.Range("I34").Value = 3367.375
$ws.Range("K34").Value = 3367.375
$ws.Range("M34").Value = -3165.375
$ws.Range("H58").Value = 3529
$ws.Range("I58").Value = 2238
$ws.Range("K58").Value = 2238
$ws.Range("M58").Value = -2035
$ws.Range("H136").Value = 3529
$ws.Range("I136").Value = 2238
$ws.Range("K136").Value = 6714
$ws.Range("M136").Value = -4164
$ws.Range("H138").Value = 40000
$ws.Range("I138").Value = 40000
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 40000
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -34860
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("H140").Value = 75000
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 119599.8
$ws.Range("J141").Value = 119599.8
$ws.Range("L141").Value = 119599.8
$ws.Range("N141").Value = -129959.8

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4343.375
$ws.Range("J80").Value = 4107
$ws.Range("L80").Value = 12321
$ws.Range("N80").Value = -14193
$ws.Range("H83").Value = 4343.375
$ws.Range("J83").Value = 4107
$ws.Range("L83").Value = 36963
$ws.Range("N83").Value = -46323
$ws.Range("H126").Value = 725
$ws.Range("I126").Value = 1266.6666
$ws.Range("K126").Value = 3799.9998
$ws.Range("M126").Value = 1140.0002

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 29998.5
$ws.Range("J57").Value = 29998.5
$ws.Range("L57").Value = 29998.5
$ws.Range("N57").Value = -31638.5
$ws.Range("H97").Value = 1713.7142
$ws.Range("I97").Value = 1754.091
$ws.Range("J97").Value = 1565.6666
$ws.Range("K97").Value = 1754.091
$ws.Range("L97").Value = 1565.6666
$ws.Range("M97").Value = -1258.091
$ws.Range("N97").Value = -2557.6666
$ws.Range("H113").Value = 3004.8948
$ws.Range("I113").Value = 1859.4
$ws.Range("K113").Value = 1859.4
$ws.Range("M113").Value = 310.5999999999999

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5616.381
$ws.Range("I16").Value = 4246.625
$ws.Range("J16").Value = 9999.6
$ws.Range("K16").Value = 4246.625
$ws.Range("L16").Value = 9999.6
$ws.Range("M16").Value = -4076.625
$ws.Range("N16").Value = -10339.6
$ws.Range("H55").Value = 384.55
$ws.Range("I55").Value = 330.6
$ws.Range("K55").Value = 330.6
$ws.Range("M55").Value = -157.6
$ws.Range("H68").Value = 2665.9167
$ws.Range("J68").Value = 2748.9
$ws.Range("L68").Value = 2748.9
$ws.Range("N68").Value = -4246.9
$ws.Range("H71").Value = 2665.9167
$ws.Range("J71").Value = 2748.9
$ws.Range("L71").Value = 13744.5
$ws.Range("N71").Value = -21232.5
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 5526
$ws.Range("I132").Value = 5080.125
$ws.Range("J132").Value = 6239.4
$ws.Range("K132").Value = 15240.375
$ws.Range("L132").Value = 18718.2
$ws.Range("M132").Value = -12710.375
$ws.Range("N132").Value = -23778.2

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 927.5333000000001
$ws.Range("J107").Value = 2800.6667
$ws.Range("L107").Value = 8402.000100000001
$ws.Range("N107").Value = -12242.0001
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H112").Value = 42387
$ws.Range("J112").Value = 42387
$ws.Range("L112").Value = 42387
$ws.Range("N112").Value = -45341
$ws.Range("H113").Value = 1413.2424
$ws.Range("I113").Value = 1231.65
$ws.Range("J113").Value = 1692.6154
$ws.Range("K113").Value = 3694.95
$ws.Range("L113").Value = 5077.8462
$ws.Range("M113").Value = -1524.95
$ws.Range("N113").Value = -9417.8462
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H122").Value = 1721.2
$ws.Range("I122").Value = 1721.2
$ws.Range("K122").Value = 5163.6
$ws.Range("M122").Value = -2713.6

Write-Host "Seraph_Profits sheets updated."
